$wb = $excel.ActiveWorkbook

# --- Overview sheet: widen columns E and F (zh-cn / de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
# Status text updated (shared string also used by de-de sheet)
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
# Latest Handback DateTime refreshed
$wsZh.Range("K2").Value = "2016-09-05 16:58:29"
# Error Detail cleared (handback version is now in sync)
$wsZh.Range("P2").Value = ""
# Column widths: Status (C) widened, Error Detail (P) narrowed
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(16).ColumnWidth = 12.833333333333332

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
# Status text updated
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
# Latest Handback DateTime refreshed
$wsDe.Range("K2").Value = "2016-09-05 16:58:37"
# Error Detail cleared
$wsDe.Range("P2").Value = ""
# Column widths: Status (C) widened, Error Detail (P) narrowed
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(16).ColumnWidth = 12.833333333333332
